$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.852.28"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "2.091.79"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "2.389.11"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "2.092.57"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "37.813.04"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.140"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.45%  "
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("E32").Value = "  +3.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("E40").Value = "  +7.16%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0213"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").Value = "1.451.84"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  +3.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("E48").Value = "  -8.49%  "
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "2.285.04"
$ws.Range("E51").Value = "  +0.99%  "
